$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $rng = $ws.Range($rangeAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '36.969.53'
$ws.Range("E2").Value = '  -0.42%  '
Set-TextValue "D3" '2.053.17'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("E4").Value = '  +0.00%  '
Set-TextValue "D5" '246.98'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("E6").Value = '  -2.16%  '
Set-TextValue "D7" '57.86'
$ws.Range("E7").Value = '  -3.60%  '
$ws.Range("E8").Value = '  +0.07%  '
Set-TextValue "D9" '0.373'
$ws.Range("E9").Value = '  -4.02%  '
Set-TextValue "D10" '0.0783'
$ws.Range("E10").Value = '  -1.49%  '
$ws.Range("E11").Value = '  +2.36%  '
Set-TextValue "D12" '15.29'
$ws.Range("E12").Value = '  -5.03%  '
Set-TextValue "D13" '0.875'
$ws.Range("E13").Value = '  +4.75%  '
Set-TextValue "D14" '2.355.79'
$ws.Range("E14").Value = '  +0.30%  '
Set-TextValue "D15" '5.65'
$ws.Range("E15").Value = '  -2.94%  '
Set-TextValue "D16" '2.087.73'
$ws.Range("E16").Value = '  +1.86%  '
Set-TextValue "D17" '18.02'
$ws.Range("E17").Value = '  -1.42%  '
Set-TextValue "D18" '36.913.75'
$ws.Range("E18").Value = '  -0.60%  '
Set-TextValue "D19" '73.80'
$ws.Range("E19").Value = '  -3.40%  '
Set-TextValue "D20" '0.0₃0895'
$ws.Range("E20").Value = '  -1.25%  '
Set-TextValue "D21" '5.42'
$ws.Range("E21").Value = '  -0.21%  '
Set-TextValue "D22" '236.15'
$ws.Range("E22").Value = '  -1.11%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("E24").Value = '  +1.69%  '
Set-TextValue "D25" '10.33'
$ws.Range("E25").Value = '  +9.22%  '
Set-TextValue "D26" '2.24'
$ws.Range("E26").Value = '  +1.16%  '
Set-TextValue "D27" '168.94'
$ws.Range("E27").Value = '  -0.08%  '
Set-TextValue "D28" '19.99'
$ws.Range("E28").Value = '  -1.24%  '
Set-TextValue "D29" '5.58'
$ws.Range("E29").Value = '  +15.19%  '
$ws.Range("E30").Value = '  -1.86%  '
Set-TextValue "D31" '1.11'
$ws.Range("E31").Value = '  -3.03%  '
Set-TextValue "D32" '4.74'
$ws.Range("E32").Value = '  +2.09%  '
Set-TextValue "D33" '0.0616'
$ws.Range("E33").Value = '  -2.47%  '
Set-TextValue "D34" '2.36'
$ws.Range("E34").Value = '  +5.20%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  +4.41%  '
Set-TextValue "D37" '0.0818'
$ws.Range("E37").Value = '  -7.68%  '
Set-TextValue "D38" '1.32'
$ws.Range("E38").Value = '  -1.72%  '
Set-TextValue "D39" '5.18'
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("E40").Value = '  -5.80%  '
Set-TextValue "D41" '0.0224'
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("E42").Value = '  +0.80%  '
Set-TextValue "D43" '0.0952'
$ws.Range("E43").Value = '  -12.16%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D44" '17.01'
$ws.Range("E44").Value = '  -3.62%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D45" '96.96'
$ws.Range("E45").Value = '  -0.33%  '
Set-TextValue "D46" '1.307.12'
$ws.Range("E46").Value = '  +0.87%  '
$ws.Range("E47").Value = '  -5.29%  '
$ws.Range("E48").Value = '  +0.01%  '
Set-TextValue "D49" '6.76'
$ws.Range("E49").Value = '  -1.30%  '
Set-TextValue "D50" '2.240.48'
$ws.Range("E50").Value = '  +0.05%  '
Set-TextValue "D51" '45.21'
$ws.Range("E51").Value = '  +1.69%  '
